$d = $word.ActiveDocument

# --- Insert "Üçüncü paragraf" right after the 4th paragraph
#     ("Burası hikayemin 2. paragrafı") ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "Üçüncü paragraf"

# --- Insert the new list paragraph right after it ---
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Baba ve oğulun hikayesi. Burada bize verilen görev hikayenin bu kısmını tamamlamak."

# --- Give it the "List Paragraph" style (saved with the Turkish
#     style id "ListeParagraf", matching the localized Word UI) ---
$p6.Style = "ListeParagraf"
$s = $d.Styles.Item("ListeParagraf")
$s.NameLocal = "List Paragraph"
$s.BaseStyle = "Normal"
$s.Priority = 34
$s.QuickStyle = $true
$s.ParagraphFormat.LeftIndent = 36
$s.NoSpaceBetweenParagraphsOfSameStyle = $true

# --- Turn it into a single-level bulleted list item (numId=1, ilvl=0) ---
$p6.Range.ListFormat.ApplyBulletDefault()
